# Applies the Polish translation updates described by the diff.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "Travel checklist: here's what you need" `
             "Podróżna lista kontrolna: oto, czego Państwo potrzebują"

Replace-Text "Here’s a checklist of the necessary items for your trip: " `
             "Oto lista przedmiotów niezbędnych podczas podróży: "

Replace-Text "Passport " "Paszport "

Replace-Text "For travellers from yellow fever endemic countries, follow the requirements set by your country. Vaccination should be done no less than 14 days prior to the journey. " `
             "For travellers from yellow fever endemic countries, follow the requirements set by your country. Szczepienie należy wykonać nie później niż 14 dni przed podróżą. "

Replace-Text "A digital or printed copy of the travel itinerary" `
             "Cyfrowa lub wydrukowana kopia planu podróży"

Replace-Text "Smart casual attire for the conference" `
             "Elegancki, swobodny strój na konferencję"

Replace-Text "Black tie attire for the Gala dinner" `
             "Czarny strój na uroczystą kolację"
